$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "NATIONAL WATCH"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = "NATIONAL WATCH"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'"
$ws.Range("O2").Value = "'"
$ws.Range("P2").Value = "'"
$ws.Range("Q2").Value = "'"
$ws.Range("R2").Value = "'"
$ws.Range("S2").Value = "'"
$ws.Range("T2").Value = "'"

# Row 3
$ws.Range("C3").Value = "GALLET"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "'"
$ws.Range("F3").Value = "GALLET"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'"
$ws.Range("O3").Value = "'"
$ws.Range("P3").Value = "'"
$ws.Range("Q3").Value = "'"
$ws.Range("R3").Value = "'"
$ws.Range("S3").Value = "'"
$ws.Range("T3").Value = "'"

# Row 4
$ws.Range("C4").Value = "HELVETIA"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "'"
$ws.Range("F4").Value = "HELVETIA"
$ws.Range("L4").Value = "'"
$ws.Range("M4").Value = "'"
$ws.Range("O4").Value = "'"
$ws.Range("P4").Value = "'"
$ws.Range("Q4").Value = "'"
$ws.Range("R4").Value = "'"
$ws.Range("S4").Value = "'"
$ws.Range("T4").Value = "'"

# Row 5
$ws.Range("C5").Value = "BELL & ROSS REF. BR 01-97 PVD STEEL LIMITED EDITION`nBell & Ross"
$ws.Range("D5").Value = "'"
$ws.Range("E5").Value = "BR"
$ws.Range("F5").Value = "BELL & ROSS REF. BR 01-97 PVD STEEL LIMITED EDITION`nBell & Ross"
$ws.Range("L5").Value = "'"
$ws.Range("M5").Value = "'"
$ws.Range("O5").Value = "'"
$ws.Range("P5").Value = "'"
$ws.Range("Q5").Value = "'"
$ws.Range("R5").Value = "'"
$ws.Range("S5").Value = "'"
$ws.Range("T5").Value = "'"

# Row 6
$ws.Range("C6").Value = "CHRONOSWISS KLASSIK REF. CH 7443 CHRONOGRAPH STEEL `nChronoswiss"
$ws.Range("D6").Value = "7443"
$ws.Range("E6").Value = "CH"
$ws.Range("F6").Value = "CHRONOSWISS KLASSIK REF. CH 7443 CHRONOGRAPH STEEL `nChronoswiss"
$ws.Range("L6").Value = "'"
$ws.Range("M6").Value = "'"
$ws.Range("O6").Value = "'"
$ws.Range("P6").Value = "'"
$ws.Range("Q6").Value = "'"
$ws.Range("R6").Value = "'"
$ws.Range("S6").Value = "'"
$ws.Range("T6").Value = "'"

# Row 7
$ws.Range("C7").Value = "VACHERON & CONSTANTIN REF. 33093 YELLOW GOLD`nVacheron & Constantin"
$ws.Range("D7").Value = "'"
$ws.Range("E7").Value = "33093"
$ws.Range("F7").Value = "VACHERON & CONSTANTIN REF. 33093 YELLOW GOLD`nVacheron & Constantin"
$ws.Range("L7").Value = "'"
$ws.Range("M7").Value = "'"
$ws.Range("O7").Value = "'"
$ws.Range("P7").Value = "'"
$ws.Range("Q7").Value = "'"
$ws.Range("R7").Value = "'"
$ws.Range("S7").Value = "'"
$ws.Range("T7").Value = "'"

# Row 8
$ws.Range("C8").Value = "ORBITA 6 ROTOR WATCH `nWINDING CABINET WOOD `nOrbita"
$ws.Range("D8").Value = "'"
$ws.Range("E8").Value = "'"
$ws.Range("F8").Value = "ORBITA 6 ROTOR WATCH `nWINDING CABINET WOOD `nOrbita"
$ws.Range("L8").Value = "'"
$ws.Range("M8").Value = "'"
$ws.Range("O8").Value = "'"
$ws.Range("P8").Value = "'"
$ws.Range("Q8").Value = "'"
$ws.Range("R8").Value = "'"
$ws.Range("S8").Value = "'"
$ws.Range("T8").Value = "'"

# Row 9
$ws.Range("C9").Value = "ORBITA 12 ROTOR WATCH `nWINDING CABINET WOOD `nOrbita"
$ws.Range("D9").Value = "'"
$ws.Range("E9").Value = "'"
$ws.Range("F9").Value = "ORBITA 12 ROTOR WATCH `nWINDING CABINET WOOD `nOrbita"
$ws.Range("L9").Value = "'"
$ws.Range("M9").Value = "'"
$ws.Range("O9").Value = "'"
$ws.Range("P9").Value = "'"
$ws.Range("Q9").Value = "'"
$ws.Range("R9").Value = "'"
$ws.Range("S9").Value = "'"
$ws.Range("T9").Value = "'"

# Row 10
$ws.Range("C10").Value = "JAEGER-LECOULTRE DESK CLOCK BRASS `nJaeger-LeCoultre"
$ws.Range("D10").Value = "'"
$ws.Range("E10").Value = "'"
$ws.Range("F10").Value = "JAEGER-LECOULTRE DESK CLOCK BRASS `nJaeger-LeCoultre"
$ws.Range("L10").Value = "'"
$ws.Range("M10").Value = "'"
$ws.Range("O10").Value = "'"
$ws.Range("P10").Value = "'"
$ws.Range("Q10").Value = "'"
$ws.Range("R10").Value = "'"
$ws.Range("S10").Value = "'"
$ws.Range("T10").Value = "'"

# Row 11
$ws.Range("C11").Value = "OFFICINE PANERAI LUMINOR MARINA PAM 111 STEEL`nOfficine Panerai"
$ws.Range("D11").Value = "'"
$ws.Range("E11").Value = "'"
$ws.Range("F11").Value = "OFFICINE PANERAI LUMINOR MARINA PAM 111 STEEL`nOfficine Panerai"
$ws.Range("L11").Value = "'"
$ws.Range("M11").Value = "'"
$ws.Range("O11").Value = "'"
$ws.Range("P11").Value = "'"
$ws.Range("Q11").Value = "'"
$ws.Range("R11").Value = "'"
$ws.Range("S11").Value = "'"
$ws.Range("T11").Value = "'"

# Row 12
$ws.Range("C12").Value = "PANERAI PAM 312 LUMINOR MARINA STEEL `nOfficine Panerai"
$ws.Range("D12").Value = "'"
$ws.Range("E12").Value = "'"
$ws.Range("F12").Value = "PANERAI PAM 312 LUMINOR MARINA STEEL `nOfficine Panerai"
$ws.Range("L12").Value = "'"
$ws.Range("M12").Value = "'"
$ws.Range("O12").Value = "'"
$ws.Range("P12").Value = "'"
$ws.Range("Q12").Value = "'"
$ws.Range("R12").Value = "'"
$ws.Range("S12").Value = "'"
$ws.Range("T12").Value = "'"

# Row 13
$ws.Range("C13").Value = "PANERAI PAM 082 MARINA MILITARE `"AMERIGO `nVESPUCCI`" TITANIUM`nOfficine Panerai"
$ws.Range("D13").Value = "'"
$ws.Range("E13").Value = "'"
$ws.Range("F13").Value = "PANERAI PAM 082 MARINA MILITARE `"AMERIGO `nVESPUCCI`" TITANIUM`nOfficine Panerai"
$ws.Range("L13").Value = "'"
$ws.Range("M13").Value = "'"
$ws.Range("O13").Value = "'"
$ws.Range("P13").Value = "'"
$ws.Range("Q13").Value = "'"
$ws.Range("R13").Value = "'"
$ws.Range("S13").Value = "'"
$ws.Range("T13").Value = "'"
